$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row pairs (2,6), (3,7), (4,8), (5,9) for columns
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).
# Column L (Calidad) stays fixed per row position (Especial/Primera/Segunda/Tercera).

$pairs = @(2, 3, 4, 5)
$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($topRow in $pairs) {
    $bottomRow = $topRow + 4

    foreach ($col in $cols) {
        $topCell = $ws.Range("$col$topRow")
        $bottomCell = $ws.Range("$col$bottomRow")

        $topVal = $topCell.Value()
        $bottomVal = $bottomCell.Value()

        $topCell.Value = $bottomVal
        $bottomCell.Value = $topVal
    }
}
